$d = $word.ActiveDocument

# caseConditional(M2DocEvaluator.java:1438) -> 1477 (single occurrence)
$d.Content.Find.Execute(
    "caseConditional(M2DocEvaluator.java:1438)", $true, $false, $false, $false, $false,
    $true, 1, $false, "caseConditional(M2DocEvaluator.java:1477)", 2) | Out-Null

# doSwitch(M2DocEvaluator.java:1216) -> 1239 (3 occurrences, all change identically)
$d.Content.Find.Execute(
    "doSwitch(M2DocEvaluator.java:1216)", $true, $false, $false, $false, $false,
    $true, 1, $false, "doSwitch(M2DocEvaluator.java:1239)", 2) | Out-Null

# caseBlock(M2DocEvaluator.java:1425) -> 1464 (single occurrence)
$d.Content.Find.Execute(
    "caseBlock(M2DocEvaluator.java:1425)", $true, $false, $false, $false, $false,
    $true, 1, $false, "caseBlock(M2DocEvaluator.java:1464)", 2) | Out-Null

# caseDocumentTemplate(M2DocEvaluator.java:287) -> 296 (single occurrence)
$d.Content.Find.Execute(
    "caseDocumentTemplate(M2DocEvaluator.java:287)", $true, $false, $false, $false, $false,
    $true, 1, $false, "caseDocumentTemplate(M2DocEvaluator.java:296)", 2) | Out-Null

# generate(M2DocEvaluator.java:276) -> 281 (single occurrence)
$d.Content.Find.Execute(
    "generate(M2DocEvaluator.java:276)", $true, $false, $false, $false, $false,
    $true, 1, $false, "generate(M2DocEvaluator.java:281)", 2) | Out-Null

# M2DocUtils.generate(M2DocUtils.java:694) -> 805 (single occurrence)
$d.Content.Find.Execute(
    "M2DocUtils.generate(M2DocUtils.java:694)", $true, $false, $false, $false, $false,
    $true, 1, $false, "M2DocUtils.generate(M2DocUtils.java:805)", 2) | Out-Null

# prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480) -> 511 (single occurrence)
$d.Content.Find.Execute(
    "prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)", $true, $false, $false, $false, $false,
    $true, 1, $false, "prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:511)", 2) | Out-Null

# generation(AbstractTemplatesTestSuite.java:389) -> 420 (single occurrence)
$d.Content.Find.Execute(
    "generation(AbstractTemplatesTestSuite.java:389)", $true, $false, $false, $false, $false,
    $true, 1, $false, "generation(AbstractTemplatesTestSuite.java:420)", 2) | Out-Null

# Insert a new stack frame line before the second RunAfters.evaluate occurrence only.
# Anchor on the unique sequence: "...ParentRunner$2.evaluate(ParentRunner.java:268)<TAB>at org...RunAfters.evaluate(RunAfters.java:27)"
$tab = [char]9
$nl = [char]10
$oldAnchor = "ParentRunner`$2.evaluate(ParentRunner.java:268)" + $nl + $tab + "at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)"
$newAnchor = "ParentRunner`$2.evaluate(ParentRunner.java:268)" + $nl + $tab + "at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)" + $nl + $tab + "at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)"

$d.Content.Find.Execute(
    $oldAnchor, $true, $false, $false, $false, $false,
    $true, 1, $false, $newAnchor, 2) | Out-Null
